# "Modifications de dernière minute"
#
# 1) Slide-number placeholder field text: ‹#› -> ‹N°› on the Slide Master
#    and on every Slide Layout.
# 2) Scorecard table on slide 5: "Angular2" -> "Angular" (table header cell).
# 3) Free-standing "Angular 2" label on slide 8 -> "Angular".

$p = $ppt.ActivePresentation

$oldNum = [string][char]0x2039 + [char]0x23 + [char]0x203A          # ‹#›
$newNum = [string][char]0x2039 + "N" + [char]0xB0 + [char]0x203A    # ‹N°›

function Update-SlideNumberPlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -ne 0) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldNum) {
                $tr.Text = $newNum
            }
        }
    }
}

# Slide Master
$master = $p.SlideMaster
Update-SlideNumberPlaceholder $master.Shapes

# Every Slide Layout hanging off the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-SlideNumberPlaceholder $layout.Shapes
}

# Scorecard table on slide 5 ("Angular2" header -> "Angular")
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shp = $slide5.Shapes.Item($i)
    if ($shp.HasTable -ne 0) {
        $tbl = $shp.Table
        for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
            for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                $cellTr = $tbl.Cell($r, $c).Shape.TextFrame.TextRange
                if ($cellTr.Text -eq "Angular2") {
                    $cellTr.Text = "Angular"
                }
            }
        }
    }
}

# Free-standing "Angular 2" caption on slide 8 -> "Angular"
$slide8 = $p.Slides.Item(8)
for ($i = 1; $i -le $slide8.Shapes.Count; $i++) {
    $shp = $slide8.Shapes.Item($i)
    if ($shp.HasTextFrame -ne 0) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "Angular 2") {
            $tr.Text = "Angular"
        }
    }
}
